# Update cryptos list values (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" column values are digit/period strings that Excel's COM
# layer would otherwise auto-coerce into numbers (dropping trailing
# zeros, collapsing thousands separators, etc). Force the cell to a
# text format right before writing so the literal string is preserved,
# matching the inlineStr the sheet already used. Only do this for cells
# whose D-value actually changes, so cells that keep their original
# value/style are left untouched.
function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

Set-TextCell "D2" "71.413.70"
Set-Cell "E2" "  +2.81%  "

Set-TextCell "D3" "3.629.30"
Set-Cell "E3" "  +6.63%  "

Set-Cell "E4" "  +0.27%  "

Set-TextCell "D5" "588.20"
Set-Cell "E5" "  +0.07%  "

Set-TextCell "D6" "179.73"
Set-Cell "E6" "  -1.06%  "

Set-TextCell "D7" "3.618.80"
Set-Cell "E7" "  +6.63%  "

Set-TextCell "D8" "0.614"
Set-Cell "E8" "  +2.29%  "

Set-Cell "E9" "  +0.09%  "

Set-TextCell "D10" "0.201"
Set-Cell "E10" "  -0.81%  "

Set-TextCell "D11" "0.605"
Set-Cell "E11" "  +1.78%  "

Set-TextCell "D12" "49.47"
Set-Cell "E12" "  +1.90%  "

Set-TextCell "D13" "0.0000285"
Set-Cell "E13" "  -1.09%  "

Set-TextCell "D14" "680.37"
Set-Cell "E14" "  -1.10%  "

Set-TextCell "D15" "4.215.53"
Set-Cell "E15" "  +6.58%  "

Set-TextCell "D16" "8.97"
Set-Cell "E16" "  +2.99%  "

# Rows 17 and 18 swap: WrappedBTC <-> WrappedEther
Set-Cell "B17" "WrappedEther"
Set-Cell "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D17" "3.642.36"
Set-Cell "E17" "  +6.31%  "

Set-Cell "B18" "WrappedBTC"
Set-Cell "C18" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D18" "71.604.00"
Set-Cell "E18" "  +3.00%  "

Set-Cell "E19" "  +1.71%  "

Set-TextCell "D20" "18.24"
Set-Cell "E20" "  +2.70%  "

Set-TextCell "D21" "11.59"
Set-Cell "E21" "  +1.99%  "

Set-TextCell "D22" "0.931"
Set-Cell "E22" "  +2.07%  "

Set-TextCell "D23" "5.86"
Set-Cell "E23" "  +9.29%  "

Set-TextCell "D24" "17.71"
Set-Cell "E24" "  +2.01%  "

Set-TextCell "D25" "102.94"
Set-Cell "E25" "  -0.45%  "

Set-TextCell "D26" "4.00"
Set-Cell "E26" "  +1.32%  "

Set-TextCell "D27" "2.83"
Set-Cell "E27" "  +3.30%  "

Set-TextCell "D28" "9.94"
Set-Cell "E28" "  +2.08%  "

Set-TextCell "D29" "34.94"
Set-Cell "E29" "  +2.54%  "

Set-TextCell "D30" "9.17"
Set-Cell "E30" "  +3.50%  "

Set-TextCell "D31" "7.23"
Set-Cell "E31" "  +3.50%  "

Set-TextCell "D32" "4.11"
Set-Cell "E32" "  +12.78%  "

Set-TextCell "D33" "573.20"
Set-Cell "E33" "  +2.83%  "

Set-TextCell "D34" "11.30"
Set-Cell "E34" "  +1.16%  "

Set-TextCell "D35" "0.109"
Set-Cell "E35" "  +2.06%  "

Set-TextCell "D36" "59.47"
Set-Cell "E36" "  +1.48%  "

Set-TextCell "D37" "0.999"
Set-Cell "E37" "  -0.10%  "

Set-TextCell "D38" "3.664.06"
Set-Cell "E38" "  +0.02%  "

Set-TextCell "D39" "0.142"
Set-Cell "E39" "  +0.67%  "

Set-TextCell "D40" "35.43"
Set-Cell "E40" "  -1.68%  "

Set-TextCell "D41" "0.0₃0757"
Set-Cell "E41" "  +3.00%  "

Set-TextCell "D42" "0.0468"
Set-Cell "E42" "  +8.97%  "

Set-TextCell "D43" "3.38"
Set-Cell "E43" "  +3.22%  "

Set-TextCell "D44" "2.74"
Set-Cell "E44" "  +1.99%  "

Set-TextCell "D45" "0.345"
Set-Cell "E45" "  +1.58%  "

Set-Cell "E46" "  +1.09%  "

Set-TextCell "D47" "2.79"
Set-Cell "E47" "  +4.03%  "

Set-Cell "E48" "  +2.21%  "

Set-TextCell "D49" "1.42"
Set-Cell "E49" "  +2.27%  "

Set-TextCell "D50" "1.00"
Set-Cell "E50" "  +0.12%  "

Set-TextCell "D51" "2.99"
Set-Cell "E51" "  +12.77%  "
